$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tbl_Admin_ParEnt")

# Update cell values (TPN test scenarios)
$ws.Range("F17").Value = 0.0775
$ws.Range("M17").Value = 0.0775

$ws.Range("M28").Value = 0

$ws.Range("C33").Value = 0.0714
$ws.Range("M33").Value = 0.103

$ws.Range("B34").Value = 0.3
$ws.Range("M34").Value = 0.0629

$ws.Range("M35").Value = 0.072

$ws.Range("M36").Value = 0.056

$ws.Range("M37").Value = 0.075

$ws.Range("L44").Value = 6.67
$ws.Range("L45").Value = 6.67

# Update the active selection on the sheet view to A28
$ws.Activate()
$ws.Range("A28").Select()
